# Update nets http interfaces document.
# Adds a new "pack process" endpoint block ({apiRoot}/pack/p) to the
# "Pack" sheet, mirroring the existing "unpack process" block already
# present on the "Unpack" sheet, and refreshes the workbook's active
# sheet / selection state.

$wb = $excel.ActiveWorkbook

$packSheet   = $wb.Worksheets.Item("Pack")
$unpackSheet = $wb.Worksheets.Item("Unpack")

# ---------------------------------------------------------------------
# 1. Bring in the formatting (fills/borders/fonts/alignment) of the
#    existing "unpack process" block (Unpack!A6:I9) without touching
#    its text, so the new rows look identical to the established
#    pattern used for every other "process" endpoint in this workbook.
# ---------------------------------------------------------------------
$unpackSheet.Range("A6:I9").Copy()
$packSheet.Range("A6:I9").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Fill in the new "pack process" request/response content.
# ---------------------------------------------------------------------
$packSheet.Range("C6").Value = "{apiRoot}/pack/p"
$packSheet.Range("D6").Value = "source files"
$packSheet.Range("E6").Value = "GET"
$packSheet.Range("F6").Value = "type TNetsPackProcessReq struct {`n Src  []string ``json:""src""```n Type string   ``json:""type""```n}"
$packSheet.Range("G6").Value = "Success"
$packSheet.Range("H6").Value = 200
$packSheet.Range("I6").Value = "type TNetsPackProcessResp struct {`n Done int64 ``json:""done""```n Work int64 ``json:""work""```n}"

$packSheet.Range("G7").Value = "Failure"
$packSheet.Range("H7").Value = 400
$packSheet.Range("I7").Value = "Bad Request"

$packSheet.Range("H8").Value = 422
$packSheet.Range("I8").Value = "Unprocessable Entity"

$packSheet.Range("H9").Value = 500
$packSheet.Range("I9").Value = "Failure Reasons"

# ---------------------------------------------------------------------
# 3. Extend / add the merges needed for the new block, matching the
#    merge layout already used by the Unpack sheet's process block.
# ---------------------------------------------------------------------
$packSheet.Range("A2:A9").Merge()
$packSheet.Range("B2:B9").Merge()
$packSheet.Range("C6:C9").Merge()
$packSheet.Range("D6:D9").Merge()
$packSheet.Range("E6:E9").Merge()
$packSheet.Range("F6:F9").Merge()
$packSheet.Range("G7:G9").Merge()

# Row 6 holds the wrapped Go struct text, so it needs the same taller
# row height already used for the analogous rows on the Unpack sheet.
$packSheet.Rows.Item(6).RowHeight = 85.5

# ---------------------------------------------------------------------
# 4. Refresh view state: Unpack's lingering selection moves onto the
#    new block, then Pack becomes the active sheet/tab with its
#    selection sitting on the new block too.
# ---------------------------------------------------------------------
$unpackSheet.Activate()
$unpackSheet.Range("E18:E21").Select()

$packSheet.Activate()
$packSheet.Range("L10:L15").Select()
